# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) timestamps on the first data
# row (row 2) of the "zh-cn" and "de-de" worksheets to reflect the new
# handback run times.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-12 02:46:03"
$wsZhCn.Range("H2").Value = "2016-03-12 02:46:30"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-12 02:46:06"
$wsDeDe.Range("H2").Value = "2016-03-12 02:46:35"
